$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row to append below the existing rows (row 5 is the last one):
# 2025/11/15 | 逃离鸭科夫 | 1119

# Column A ("Date") stores its entries as plain text in this workbook
# (e.g. "2025/11/12"), not as real date values. Typing a date-shaped string
# into a General-formatted cell would make Excel auto-convert it to a date
# serial number, so force the cell to Text format first, then write the
# literal string.
$dateCell = $ws.Cells.Item(6, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025/11/15"

# Re-apply the plain "Normal" base style so the earlier Text number-format
# doesn't linger, then apply the centered alignment that rows 3-5 use
# (style index s="1": horizontal/vertical center, General format).
$dateCell.Style = "Normal"
$dateCell.HorizontalAlignment = -4108  # xlCenter
$dateCell.VerticalAlignment = -4108    # xlCenter

$ws.Cells.Item(6, 2).Value = "逃离鸭科夫"
$ws.Cells.Item(6, 3).Value = 1119

$gameAndCount = $ws.Range("B6:C6")
$gameAndCount.HorizontalAlignment = -4108  # xlCenter
$gameAndCount.VerticalAlignment = -4108    # xlCenter
